# "emptied device key list"
# Clear out all the device-key data rows (A2:E18), leaving the header row
# (row 1) and the existing cell formatting (style index) in place — this
# mirrors selecting the data range and pressing Delete / Clear Contents
# in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:E18")
$dataRange.ClearContents()

# Leave the range selected, same as it would be right after clearing it.
$dataRange.Select()
